$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (General/inline-string) representation
# rather than being auto-converted to numbers by Excel when values look numeric.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.895.67'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.798.77'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.89%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.007'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '304.00'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4944'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3836'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09149'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +14.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.090'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.59'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.012'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.302'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.51'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.806.61'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.189'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001103'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06555'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.007'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.00'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.939'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.987.88'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.214'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.22'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.017.75'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.38'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.347'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.55'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1063'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.038'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.621'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.520'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06796'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.806'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02291'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2119'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.30'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -7.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.902'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6089'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.142'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.283'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5840'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.654'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.22'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.932'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.167'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06740'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.12%  '
